$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
  @(2.0441090480123867E-3, 0.1046576942724959, 2.6315870816383365E-2, 818923843.82635462, 23.826098078463101, 61.845370066492791, 0.6884271576465042, 6.8694764167520919E-2, 9.9999999572690911E-3, 5.0746472033720158E-11, 0.99150334061473511, 0.98468599961749703, 0.96903446615868649, 9.1715281122680306E-2, 0.12312944561995699, 0.17508820206655948, 2.1256678515456123E-2, 2.6613300701501686E-2),
  @(1.7375150690284137E-3, 0.10496541676836456, 2.6315789474102722E-2, 40932149.670228243, 18.766309887571968, 59.827350501386753, 1.000000018820988E-2, 2.6610097734376798E-2, 3.3355305486471564E-14, 9.9999999999766283E-3, 0.99651851479241171, 0.97294806800969846, 0.99312022183013082, 5.8708349964144334E-2, 0.16365027549747219, 8.2528664039656183E-2, 1.9607843138263741E-2, 2.6549236112445213E-2),
  @(2.5187632818234077E-3, 6.8831324853419681E-2, 9.9974585701027421E-2, 1041012019.7963969, 19.199057782624216, 68.746266970956583, 0.60195391607972648, 1.3554373795349711E-2, 1.8592919601092178E-4, 2.8977858246021342E-8, 0.99165914719301107, 0.98341885734021683, 0.99354923826659625, 9.0870480789523103E-2, 0.12812232917483732, 7.9914042045623468E-2, 2.0762247073926928E-2, 2.6530005734012874E-2),
  @(3.8069722426232336E-3, 0.18523094004241433, 2.6318056725626454E-2, 96657727.729502618, 31.026684917817725, 59.368951390179589, 1.0000000075770199E-2, 0.10090365131518088, 2.8046324163883419E-13, 9.9999999999640272E-3, 0.99344509117279267, 0.91661952744682262, 0.88432615982790008, 8.0556562357980574E-2, 0.287309359102076, 0.33840375555005137, 1.9607843137974972E-2, 2.6647497481135117E-2),
  @(8.7698902232442831E-3, 0.86399580511486929, 9.995881077473788E-2, 1082559483.0567219, 12.230932328550869, 70.551362167693284, 0.56959627414767056, 0.60703740666618455, 2.1880786865386929E-6, 6.7874104622766673E-3, 0.99461185579455313, 0.97388629316479802, 0.97501087959352273, 7.3036037429425077E-2, 0.16078734330428476, 0.15728709165857338, 2.0709893300277678E-2, 2.7405600634281477E-2),
  @(7.3264351026778645E-3, 0.24643861541745188, 2.6315889911913106E-2, 109078677.38140464, 23.19693052034668, 55.748231681964143, 1.0000018886960196E-2, 8.3637230584619873E-2, 4.603846508116658E-7, 9.9999999990016877E-3, 0.99557530586700138, 0.98326944414323225, 0.9832443088179541, 6.6184946866101421E-2, 0.12869829174546227, 0.12879493107349146, 1.9607843178885108E-2, 2.6623949238280369E-2),
  @(7.8106083261215891E-3, 5.869961277407737E-2, 3.5455462328748362E-2, 869498099.03153777, 39.085223331083384, 54.559074628585712, 0.80285495997640921, 4.7834282713436441E-2, 9.999999998044606E-3, 3.37909909116856E-14, 0.99705613162078521, 0.99612598591023427, 0.99622230889038665, 5.3985458184798663E-2, 6.1929588637969039E-2, 6.1154837899525014E-2, 2.1256907874162789E-2, 2.6538930104340355E-2),
  @(1.0293787161300128E-2, 6.2543215632686569E-2, 9.9999999454831554E-2, 1105943839.9999909, 34.224737576936299, 62.585372843992985, 0.66635472008517105, 7.9144405473099394E-2, 3.4630397152366507E-9, 8.553147623255936E-11, 0.99748717120450758, 0.99501459696272554, 0.99238595144909525, 4.9879347454309449E-2, 7.0256975533841526E-2, 8.6825468219807758E-2, 2.1040263260226788E-2, 2.6618504668151705E-2),
  @(3.895568585440214E-3, 0.38617223945573631, 9.9999999999977801E-2, 489079715.7061702, 29.434333457740525, 66.762303612276114, 7.9001729104433732E-2, 0.9999999999999778, 3.22176319993991E-14, 9.9987983014998855E-3, 0.99576792723484286, 0.9246492770850635, 0.95852967094049868, 6.4728293948671162E-2, 0.27312490857808475, 0.20262187880114618, 1.9939613420300348E-2, 2.8011204481770458E-2),
  @(3.0912767097399747E-3, 0.27786428949928771, 9.9970729895002844E-2, 61762617.799501352, 37.503376682768554, 39.887875194800991, 2.96468843504391E-2, 0.24161039050374877, 8.8736851322044502E-5, 6.1827917279074843E-5, 0.99451206993489205, 0.98430273911189758, 0.99124092134438857, 7.3709231202454359E-2, 0.12466069259883569, 9.3120824035525379E-2, 1.9675462070291562E-2, 2.6872708758745849E-2),
  @(8.5966173476283748E-3, 7.5576964837486896E-2, 4.083205714724527E-2, 110535997.81545092, 18.843475719718683, 71.952028245861996, 0.61767771816427552, 9.8468690080231028E-2, 5.3488146179637036E-10, 5.4600006759997276E-11, 0.98377993628833738, 0.98268618915063599, 0.98190427074829167, 0.12671962387312385, 0.13092239205296552, 0.13384607562118209, 2.1225779808462727E-2, 2.6593503204189375E-2),
  @(2.484422966876763E-2, 0.548688174495364, 9.9999999999454886E-2, 1451043699.9426327, 12.000000000136014, 77.205567009344776, 0.50226857233473965, 0.99999999999997524, 7.1581570957742841E-8, 2.2367537042701757E-10, 0.96577830725755187, 0.92024514422278103, 0.98066981225722893, 0.18406378192089728, 0.28099342913927139, 0.13833613362149241, 2.0819631505320322E-2, 2.801120448177051E-2),
  @(5.3544805483668996E-3, 2.6088655451939755E-2, 9.9999999999976122E-2, 245215568.69201437, 28.493957329015188, 56.471391752168415, 1.0000034065928785E-2, 1.7334773413134628E-2, 3.3331101907621433E-14, 9.9999999999757522E-3, 0.99859935377952802, 0.99482873884190126, 0.99009993549429876, 3.7237612145078626E-2, 7.1551020583341687E-2, 9.90003225279808E-2, 1.9607843278592513E-2, 2.6534500435993966E-2),
  @(3.7611854946927692E-3, 2.5565931022230447E-2, 3.5321231229687819E-2, 2505455714.6952477, 33.023716670948112, 59.464318160760158, 0.73350512826641245, 1.4414486891105648E-2, 9.9327927916984752E-3, 1.3208270856694955E-4, 0.99716412180665903, 0.99332643771301776, 0.99212247336633141, 5.2986030342039983E-2, 8.1282388400639324E-2, 8.831053938988187E-2, 2.1273982818421418E-2, 2.6532383197032189E-2),
  @(5.1199478980145652E-3, 4.7267980561571388E-2, 6.2905539435138727E-2, 246228258.5315257, 18.180706790665628, 67.337308865104603, 1.0000000000039561E-2, 1.0743085172171422E-2, 3.9290826589744681E-14, 9.9999999975049127E-3, 0.99864701909399545, 0.98413426305868923, 0.99348792197162894, 3.6598512223100878E-2, 0.12532788824478638, 8.0292946440439633E-2, 1.960784313729446E-2, 2.6526016201373685E-2),
  @(2.6357140665028162E-3, 0.22097747374057819, 9.9999999999976386E-2, 1520030069.0057569, 42.233840979344066, 67.258455347920034, 0.65485030604715355, 0.99999999999995892, 3.2205104626606272E-14, 3.2469114290095237E-14, 0.99372489872940917, 0.96158848696192278, 0.93221528541603726, 7.881846394015124E-2, 0.19500614838434299, 0.25904993232603452, 2.1252851843615587E-2, 2.8011204481770476E-2),
  @(6.7948149305836654E-3, 0.168014430025443, 2.6316124966346988E-2, 63316808.052846134, 32.821068749902935, 56.40340088678083, 1.0000012081835695E-2, 5.1687544118648891E-2, 4.4727504173871706E-10, 9.9999875984612494E-3, 0.99612548610314, 0.97570643760554632, 0.99468610933668, 6.193358344138826E-2, 0.15508264496876856, 7.253103995315939E-2, 1.9607843157620135E-2, 2.6583961549957996E-2),
  @(2.2362757669413497E-2, 0.44317941837882197, 9.9997538266905939E-2, 1831633444.395858, 13.761915608029204, 72.559171492127277, 0.55582267615613101, 0.27845077781799188, 9.5769108440663419E-10, 5.931258553957223E-3, 0.99527498280796078, 0.96689731351638464, 0.99058588898251099, 6.8394203117799895E-2, 0.18102944406581828, 9.6539991233240516E-2, 2.1036571785094701E-2, 2.6867372780509297E-2),
  @(7.9999999999977797E-2, 8.966010563456564E-2, 2.6315789473706708E-2, 50540652.538597472, 25.913858634253469, 62.398474319054834, 0.60844375324832323, 2.7337576462009381E-2, 3.220446049250329E-14, 3.2638565114095641E-14, 0.98245617229471505, 0.9313162494548981, 0.93501958607585667, 0.13178918555113725, 0.26076217716465488, 0.25363479608464984, 2.110165987808248E-2, 2.6546559914688576E-2),
  @(1.0759029324079272E-2, 0.10806440153231296, 9.9999999999865849E-2, 1442195507.9818859, 21.071899069677666, 66.224948294347939, 0.64161205898035911, 0.11940514825193088, 2.843330400548159E-12, 9.9999988405638767E-3, 0.99428979830560293, 0.99141732978265473, 0.99019570580954175, 7.5187097812411327E-2, 9.2178324540923715E-2, 9.8520308812719781E-2, 2.0973295676043606E-2, 2.6672789715633204E-2),
  @(7.6944337044300348E-3, 0.18642578593065726, 9.9999999999977801E-2, 2066751694.6130178, 15.152082157970085, 69.229555886375422, 0.62015468114133709, 0.14492588275622542, 3.2204559023339198E-14, 9.999999999789979E-3, 0.98943801322145664, 0.98056065706875795, 0.98156150579240986, 0.10225637833777387, 0.13872616732949003, 0.13510776907903643, 2.122528494333693E-2, 2.6678562379358721E-2),
  @(6.7398310494739926E-3, 0.22290286373119217, 9.9999999791558727E-2, 1070399379.9948366, 15.623289739233876, 75.123932398316541, 0.54377409760023998, 0.18253535474546018, 8.2559360702717407E-8, 2.7534271388347482E-7, 0.99162805440929291, 0.93842197202228206, 0.9921309821534936, 9.1039695379543054E-2, 0.24690534157433841, 8.8262832880218509E-2, 2.0682745671593006E-2, 2.6719952401461845E-2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 1
  $row = $data[$i]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}

# Update bound cells Z5/AA5 (change from formula 1E-8 to literal 0.01) and Z6/AA6 (10 -> 1)
$ws.Range("Z5").Value = 0.01
$ws.Range("AA5").Value = 0.01
$ws.Range("Z6").Value = 1
$ws.Range("AA6").Value = 1

# Update sheet view: remove topLeftCell freeze position and change selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F10").Select()

# Move the "Chart 3" graphic frame (on Sheet1) to its new anchor position.
# from: col 13 offset 552994, row 46 offset 32385
# to:   col 21 offset 248194, row 60 offset 108585
$shape = $ws.Shapes.Item("Chart 3")
$shape.TopLeftCell = $ws.Cells.Item(47, 14)
$shape.Left = $ws.Columns.Item(14).Left + 552994 / 914400 * 72
$shape.Top = $ws.Rows.Item(47).Top + 32385 / 914400 * 72
